$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Absent" column (H) values to reflect the consolidated report.
$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 0
$ws.Range("H12").Value = 1
$ws.Range("H15").Value = 0
